$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to reflect the new filename ("find1" -> "find1.pdf")
$ws.Range("A2").Value = "find1.pdf"

# A3 keeps its value "notfind2" (rewritten to refresh shared-string ordering)
$ws.Range("A3").Value = "notfind2"

# Move the active selection from A4 to A8
$ws.Range("A8").Select()
